$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full data table (A1:B11)
$data = @(
    @("Data Name", "Data Value "),
    @("chargeItemValue1", "UCCITest#01"),
    @("chargeItemValue2", "UCCITest#02"),
    @("chargeItemValue3", "UCCITest#03"),
    @("chargeItemValue4", "UCCITest#04"),
    @("chargeItemDescription1", "For testing"),
    @("amount1", '"500"'),
    @("amount2", '"600.50"'),
    @("searchItemValue1", "UCCITest#01"),
    @("searchItemValue2", "UC_CI_Test#02"),
    @("expectedValue1", "UC_CI_Test#02")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

$ws.Range("F16").Select()
